$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Create the new "Week 6" sheet by copying "Week 5" (keeps all formatting,
# column widths, merges, etc. identical to the other week sheets) and
# placing it after the last sheet in the workbook.
# ---------------------------------------------------------------------------
$weekFive = $wb.Worksheets.Item("Week 5")
$weekFive.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$weekSix = $wb.Worksheets.Item($wb.Worksheets.Count)
$weekSix.Name = "Week 6"

# ---------------------------------------------------------------------------
# Replace the matchup data (home team / spread / away team / abbreviations)
# with Week 6's games. Columns: C = home team, D = spread, E = away team,
# I = home team abbreviation, K = away team abbreviation.
# ---------------------------------------------------------------------------
$games = @(
    @{ Row = 2;  Home = "PANTHERS"; Spread = 2.5; Away = "BEARS";         HomeAbbr = "CAR"; AwayAbbr = "CHI" },
    @{ Row = 3;  Home = "LIONS";    Spread = 3.5; Away = "JAGUARS";       HomeAbbr = "DET"; AwayAbbr = "JAX" },
    @{ Row = 4;  Home = "VIKINGS";  Spread = 3.5; Away = "FALCONS";       HomeAbbr = "MIN"; AwayAbbr = "ATL" },
    @{ Row = 5;  Home = "TEXANS";   Spread = 0;   Away = "TITANS";        HomeAbbr = "HOU"; AwayAbbr = "TEN" },
    @{ Row = 6;  Home = "GIANTS";   Spread = 2.5; Away = "FOOTBALL TEAM"; HomeAbbr = "NYG"; AwayAbbr = "WAS" },
    @{ Row = 7;  Home = "STEELERS"; Spread = 3.5; Away = "BROWNS";        HomeAbbr = "PIT"; AwayAbbr = "CLE" },
    @{ Row = 8;  Home = "RAVENS";   Spread = 7.5; Away = "EAGLES";        HomeAbbr = "BAL"; AwayAbbr = "PHI" },
    @{ Row = 9;  Home = "COLTS";    Spread = 8;   Away = "BENGALS";       HomeAbbr = "IND"; AwayAbbr = "CIN" },
    @{ Row = 10; Home = "BRONCOS";  Spread = 0;   Away = "PATRIOTS";      HomeAbbr = "DEN"; AwayAbbr = "NE"  },
    @{ Row = 11; Home = "DOLPHINS"; Spread = 8.5; Away = "JETS";          HomeAbbr = "MIA"; AwayAbbr = "NYJ" },
    @{ Row = 12; Home = "PACKERS";  Spread = 2.5; Away = "BUCCANEERS";    HomeAbbr = "GB";  AwayAbbr = "TB"  },
    @{ Row = 13; Home = "RAMS";     Spread = 3.5; Away = "49ERS";         HomeAbbr = "LAR"; AwayAbbr = "SF"  },
    @{ Row = 14; Home = "CHIEFS";   Spread = 0;   Away = "BILLS";         HomeAbbr = "KC";  AwayAbbr = "BUF" }
)

foreach ($game in $games) {
    $r = $game.Row
    $weekSix.Range("C$r").Value = $game.Home
    $weekSix.Range("D$r").Value = $game.Spread
    $weekSix.Range("E$r").Value = $game.Away
    $weekSix.Range("I$r").Value = $game.HomeAbbr
    $weekSix.Range("K$r").Value = $game.AwayAbbr
}

# Week 5 has 14 games (last one on row 15); Week 6 only has 13, so row 15's
# game-data cells must be cleared out on the new sheet.
$weekSix.Range("C15").ClearContents()
$weekSix.Range("D15").ClearContents()
$weekSix.Range("E15").ClearContents()
$weekSix.Range("I15").ClearContents()
$weekSix.Range("K15").ClearContents()

# ---------------------------------------------------------------------------
# The thick bottom-border formatting (columns N/O) that marks the end of the
# game list needs to move up one row to match the 13-game week: row 13 takes
# on the "last row" look, row 15 reverts to the plain blank-row look.
# ---------------------------------------------------------------------------
$weekSix.Range("N15:O15").Copy($weekSix.Range("N13"))
$weekSix.Range("N17:O17").Copy($weekSix.Range("N15"))

$weekSix.Range("A1").Select() | Out-Null
